$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet (the "name" attribute in workbook.xml, not just the tab caption)
$ws.Name = "BrassA"

# 2. Tiny last-digit value corrections in existing rows.
#    Values are assigned as strings so the engine stores the exact decimal
#    text instead of re-deriving a 17-significant-digit expansion of the
#    double (which would introduce a spurious trailing digit).
$ws.Range("C13").Value = "1.010127284223264"
$ws.Range("G13").Value = "1.010127284223264"
$ws.Range("N13").Value = "0.9966567526695027"
$ws.Range("K15").Value = "0.9261600377335842"

# 3. Append a new row 16 with data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = "1.91528448868924"
$ws.Range("D16").Value = "1.130738473172443"
$ws.Range("E16").Value = "1.183960349202367"
$ws.Range("F16").Value = "0.9949748176569884"
$ws.Range("G16").Value = "1.91528448868924"
$ws.Range("H16").Value = "1.130738473172443"
$ws.Range("I16").Value = "0.9930788932318985"
$ws.Range("J16").Value = "0.6005972044409383"
$ws.Range("K16").Value = "1.17902922840198"
$ws.Range("L16").Value = "0.9894688374312056"
$ws.Range("M16").Value = "1.91528448868924"
$ws.Range("N16").Value = "1.157349411187405"
$ws.Range("O16").Value = "1.30623953218026"
$ws.Range("P16").Value = "1.123391536528383"

# Apply the same formatting used for the other "A" column index cells
# (bold, centered/top-aligned, thin box border) to the new A16 cell by
# copying the format from A15 (avoids minting a redundant/unused style
# record, unlike setting Font/Alignment/Borders piecemeal).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
